# Fixed PB and EB DSL issues
#
# - D2 / D3 get the "A" automatable flag (already present on D4).
# - G3/H3 (uuid test) and G4/H4 (oeminfo test) DSL steps gain a
#   "validate4" step (validating a new SystemProperties check) instead of
#   the old "CheckGenericValues(uuid|oeminfo)" call.
# - Row heights for rows 3 and 4 grow to fit the extra DSL line.
# - Selection moves from C1 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- D column: mark rows 2 and 3 as Automatable ("A"), matching D4 ---
$ws.Range("D2").Value = "A"
$ws.Range("D3").Value = "A"

# NOTE: new shared-string entries are appended in first-write order, and the
# target workbook expects them interleaved as: oeminfo-steps, oeminfo-validate,
# uuid-validate, uuid-steps. Write the cells in that order so the resulting
# shared-string table lines up with the target file.

# --- Row 4 (Testcase 3, "Generic ActiveX Object OEM Information String") ---
$ws.Range("G4").Value = "wait(3);`nvalidate1;`nlink_Click(generic_test_link);`nvalidate2;`nlink_Click(activex_link);`nvalidate3;`nwait(1);`nSelectTestToRun(VT056_1357_string);`nwait(1);`nClickRunTest(runtest_top_xpath);`nvalidate4;"

$ws.Range("H4").Value = "validate1`n{`nvalidate_PageTitle=Pocket Browser Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Module Index Page`n};`nvalidate3`n{`nvalidate_PageTitle=PB and RE2.2 Semi Auto Frame Work : Generic`n};`nvalidate4`n{`nvalidate_SystemProperties=oeminfo`n};"

# --- Row 3 (Testcase 2, "Generic ActiveX Object UUID") ---
$ws.Range("H3").Value = "validate1`n{`nvalidate_PageTitle=Pocket Browser Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Module Index Page`n};`nvalidate3`n{`nvalidate_PageTitle=PB and RE2.2 Semi Auto Frame Work : Generic`n};`nvalidate4`n{`nvalidate_SystemProperties=uuid`n};`n"

$ws.Range("G3").Value = "wait(3);`nvalidate1;`nlink_Click(generic_test_link);`nvalidate2;`nlink_Click(activex_link);`nvalidate3;`nwait(1);`nSelectTestToRun(VT056_0779_string);`nwait(1);`nClickRunTest(runtest_top_xpath);`nvalidate4;"

# --- Row heights grow to fit the extra validate4 DSL line ---
$ws.Rows.Item(3).RowHeight = 203.25
$ws.Rows.Item(4).RowHeight = 192

# --- Selection moves from C1 to A2 ---
$ws.Range("A2").Select()
